$d = $word.ActiveDocument

# Find the last paragraph in the body (the "Invincible mode..." bullet)
$paragraphs = $d.Paragraphs
$lastPara = $paragraphs.Last

# Insert a new paragraph right after it
$newRange = $lastPara.Range.InsertParagraphAfter()

# The newly created paragraph is now the (new) last paragraph
$newPara = $d.Paragraphs.Last

# Match the bullet-list formatting of the paragraph it follows
$newPara.Style = $lastPara.Style
$newPara.Range.ListFormat.ListTemplate = $lastPara.Range.ListFormat.ListTemplate
$newPara.Range.ListFormat.ListLevelNumber = $lastPara.Range.ListFormat.ListLevelNumber

# Set the text of the new paragraph
$newPara.Range.Text = "Added level 2"
